$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("summary")
$wsModelFit = $wb.Worksheets.Item("model_fit")

# summary sheet updates (rows 2-16)
$wsSummary.Range("D2").Value = 964
$wsSummary.Range("E2").Value = 83.82
$wsSummary.Range("F2").Value = -2.06
$wsSummary.Range("H2").Value = 1.01
$wsSummary.Range("I2").Value = 0.14
$wsSummary.Range("J2").Value = 0.29
$wsSummary.Range("K2").Value = 0.03
$wsSummary.Range("L2").Value = 1.14
$wsSummary.Range("D3").Value = 959
$wsSummary.Range("E3").Value = 79.56
$wsSummary.Range("F3").Value = -1.71
$wsSummary.Range("H3").Value = 1
$wsSummary.Range("I3").Value = -0.06
$wsSummary.Range("J3").Value = 0.33
$wsSummary.Range("L3").Value = 1.27
$wsSummary.Range("D4").Value = 966
$wsSummary.Range("E4").Value = 76.6
$wsSummary.Range("F4").Value = -1.5
$wsSummary.Range("H4").Value = 0.98
$wsSummary.Range("I4").Value = -0.53
$wsSummary.Range("J4").Value = 0.37
$wsSummary.Range("L4").Value = 1.34
$wsSummary.Range("D5").Value = 959
$wsSummary.Range("E5").Value = 73.51
$wsSummary.Range("F5").Value = -1.29
$wsSummary.Range("I5").Value = -0.4
$wsSummary.Range("J5").Value = 0.38
$wsSummary.Range("L5").Value = 1.3
$wsSummary.Range("D6").Value = 953
$wsSummary.Range("E6").Value = 65.27
$wsSummary.Range("F6").Value = -0.8
$wsSummary.Range("I6").Value = 0.53
$wsSummary.Range("D7").Value = 946
$wsSummary.Range("E7").Value = 62.16
$wsSummary.Range("F7").Value = -0.64
$wsSummary.Range("H7").Value = 0.99
$wsSummary.Range("I7").Value = -0.21
$wsSummary.Range("J7").Value = 0.39
$wsSummary.Range("K7").Value = 0.03
$wsSummary.Range("D8").Value = 941
$wsSummary.Range("E8").Value = 57.49
$wsSummary.Range("F8").Value = -0.39
$wsSummary.Range("G8").Value = 0.07
$wsSummary.Range("H8").Value = 1.01
$wsSummary.Range("I8").Value = 0.23
$wsSummary.Range("J8").Value = 0.38
$wsSummary.Range("K8").Value = 0.04
$wsSummary.Range("L8").Value = 1.19
$wsSummary.Range("D9").Value = 926
$wsSummary.Range("E9").Value = 52.92
$wsSummary.Range("F9").Value = -0.14
$wsSummary.Range("G9").Value = 0.07
$wsSummary.Range("I9").Value = 0.47
$wsSummary.Range("J9").Value = 0.38
$wsSummary.Range("L9").Value = 1.14
$wsSummary.Range("D10").Value = 924
$wsSummary.Range("E10").Value = 46.65
$wsSummary.Range("F10").Value = 0.19
$wsSummary.Range("G10").Value = 0.07
$wsSummary.Range("H10").Value = 1.03
$wsSummary.Range("I10").Value = 1.07
$wsSummary.Range("K10").Value = 0.03
$wsSummary.Range("L10").Value = 1.05
$wsSummary.Range("D11").Value = 895
$wsSummary.Range("E11").Value = 40.78
$wsSummary.Range("F11").Value = 0.5
$wsSummary.Range("I11").Value = 0.01
$wsSummary.Range("J11").Value = 0.4
$wsSummary.Range("K11").Value = 0.02
$wsSummary.Range("L11").Value = 1.22
$wsSummary.Range("D12").Value = 842
$wsSummary.Range("E12").Value = 35.27
$wsSummary.Range("F12").Value = 0.79
$wsSummary.Range("H12").Value = 0.97
$wsSummary.Range("I12").Value = -0.88
$wsSummary.Range("J12").Value = 0.41
$wsSummary.Range("L12").Value = 1.4
$wsSummary.Range("D13").Value = 785
$wsSummary.Range("E13").Value = 31.85
$wsSummary.Range("F13").Value = 0.97
$wsSummary.Range("H13").Value = 1.02
$wsSummary.Range("I13").Value = 0.65
$wsSummary.Range("K13").Value = 0.03
$wsSummary.Range("L13").Value = 1.09
$wsSummary.Range("D14").Value = 679
$wsSummary.Range("E14").Value = 24.01
$wsSummary.Range("F14").Value = 1.45
$wsSummary.Range("H14").Value = 1.02
$wsSummary.Range("I14").Value = 0.34
$wsSummary.Range("J14").Value = 0.34
$wsSummary.Range("L14").Value = 1.15
$wsSummary.Range("D15").Value = 552
$wsSummary.Range("E15").Value = 20.65
$wsSummary.Range("F15").Value = 1.65
$wsSummary.Range("H15").Value = 0.95
$wsSummary.Range("I15").Value = -0.85
$wsSummary.Range("J15").Value = 0.39
$wsSummary.Range("K15").Value = 0.05
$wsSummary.Range("L15").Value = 1.63
$wsSummary.Range("D16").Value = 345
$wsSummary.Range("E16").Value = 16.23
$wsSummary.Range("F16").Value = 2.02
$wsSummary.Range("G16").Value = 0.16
$wsSummary.Range("H16").Value = 1
$wsSummary.Range("I16").Value = 0.08
$wsSummary.Range("J16").Value = 0.3
$wsSummary.Range("K16").Value = 0.04
$wsSummary.Range("L16").Value = 1.14

# model_fit sheet updates (rows 2-3)
$wsModelFit.Range("D2").Value = 13914
$wsModelFit.Range("E2").Value = 13946
$wsModelFit.Range("F2").Value = 14024
$wsModelFit.Range("G2").Value = 0.738
$wsModelFit.Range("H2").Value = 0.665
$wsModelFit.Range("D3").Value = 13902
$wsModelFit.Range("E3").Value = 13962
$wsModelFit.Range("F3").Value = 14110
$wsModelFit.Range("G3").Value = 0.739
$wsModelFit.Range("H3").Value = 0.666
